$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.062.67"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "1.761.36"
$ws.Range("E3").Value = "  -1.20%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").Value = "'334.77"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("D6").Value = "'0.9976"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").Value = "'0.3882"
$ws.Range("E7").Value = "  +1.77%  "

$ws.Range("D8").Value = "'0.3395"
$ws.Range("E8").Value = "  -1.36%  "

$ws.Range("E9").Value = "  -3.41%  "

$ws.Range("D10").Value = "'1.124"
$ws.Range("E10").Value = "  -2.36%  "

$ws.Range("D11").Value = "'0.07210"
$ws.Range("E11").Value = "  -2.24%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'22.44"
$ws.Range("E12").Value = "  -3.00%  "

$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").Value = "'0.9990"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("E14").Value = "  -4.55%  "

$ws.Range("D15").Value = "1.754.27"
$ws.Range("E15").Value = "  -1.73%  "

$ws.Range("D16").Value = "'7.050"
$ws.Range("E16").Value = "  -4.05%  "

$ws.Range("D17").Value = "'0.00001058"
$ws.Range("E17").Value = "  -1.36%  "

$ws.Range("D18").Value = "'0.06605"
$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("D19").Value = "'80.59"
$ws.Range("E19").Value = "  -2.09%  "

$ws.Range("D20").Value = "'0.9973"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("E21").Value = "  -3.43%  "

$ws.Range("D22").Value = "'6.205"
$ws.Range("E22").Value = "  -3.88%  "

$ws.Range("D23").Value = "28.047.67"
$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("D24").Value = "'11.63"
$ws.Range("E24").Value = "  -3.51%  "

$ws.Range("D25").Value = "'2.390"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("D26").Value = "'154.37"
$ws.Range("E26").Value = "  +0.23%  "

$ws.Range("D27").Value = "'19.91"
$ws.Range("E27").Value = "  -3.88%  "

$ws.Range("D28").Value = "'2.309"
$ws.Range("E28").Value = "  -4.39%  "

$ws.Range("D29").Value = "1.951.42"
$ws.Range("E29").Value = "  -1.78%  "

$ws.Range("D30").Value = "'1.279"
$ws.Range("E30").Value = "  -11.98%  "

$ws.Range("D31").Value = "'129.32"
$ws.Range("E31").Value = "  -5.10%  "

$ws.Range("D32").Value = "'4.063"
$ws.Range("E32").Value = "  +3.15%  "

$ws.Range("D33").Value = "'5.837"
$ws.Range("E33").Value = "  -4.60%  "

$ws.Range("D34").Value = "'0.08705"
$ws.Range("E34").Value = "  -1.87%  "

$ws.Range("D35").Value = "'12.06"
$ws.Range("E35").Value = "  -5.57%  "

$ws.Range("D36").Value = "'0.02282"
$ws.Range("E36").Value = "  -6.45%  "

$ws.Range("D37").Value = "'5.138"
$ws.Range("E37").Value = "  -3.67%  "

$ws.Range("D38").Value = "'0.06152"
$ws.Range("E38").Value = "  -3.29%  "

$ws.Range("D39").Value = "'0.6496"
$ws.Range("E39").Value = "  -5.23%  "

$ws.Range("D40").Value = "'0.2109"
$ws.Range("E40").Value = "  -3.02%  "

$ws.Range("D41").Value = "'1.496"
$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("E42").Value = "  -3.38%  "

$ws.Range("D43").Value = "'0.9970"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").Value = "'7.870"
$ws.Range("E44").Value = "  -5.14%  "

$ws.Range("D45").Value = "'13.74"
$ws.Range("E45").Value = "  -2.81%  "

$ws.Range("D46").Value = "'3.828"
$ws.Range("E46").Value = "  -1.20%  "

$ws.Range("D47").Value = "'0.6000"
$ws.Range("E47").Value = "  -4.66%  "

$ws.Range("E48").Value = "  -5.00%  "

$ws.Range("E49").Value = "  -5.13%  "

$ws.Range("D50").Value = "'0.07004"
$ws.Range("E50").Value = "  -5.88%  "

$ws.Range("D51").Value = "'1.154"
$ws.Range("E51").Value = "  -4.19%  "
